# Append the July 2021 daily subscription rows (new shared-string dates +
# numeric totals per original-term column) to the end of the existing table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Date (col A, stored as text), Total, 3m, 6m, 9m, 12m, 18m, 2y, 5y, 10y+
$data = @(
    @("01-07-2021", 623, 0, 0, 0, 84, 88, 237, 33, 180),
    @("02-07-2021", 562, 0, 0, 0, 0, 58, 179, 254, 71),
    @("05-07-2021", 3, 0, 0, 0, 0, 0, 0, 0, 3),
    @("06-07-2021", 654, 0, 168, 49, 0, 0, 258, 104, 76),
    @("07-07-2021", 174, 0, 0, 0, 0, 0, 51, 66, 57),
    @("08-07-2021", 960, 192, 450, 30, 45, 20, 102, 87, 34),
    @("09-07-2021", 469, 244, 48, 0, 37, 41, 38, 24, 37),
    @("12-07-2021", 664, 0, 0, 64, 145, 40, 183, 191, 42),
    @("13-07-2021", 711, 192, 48, 32, 44, 15, 181, 162, 37),
    @("14-07-2021", 961, 206, 59, 81, 226, 16, 97, 189, 87),
    @("15-07-2021", 1059, 294, 392, 32, 22, 0, 64, 134, 121),
    @("19-07-2021", 354, 0, 0, 62, 62, 42, 72, 42, 75),
    @("20-07-2021", 696, 0, 122, 0, 200, 58, 223, 56, 37),
    @("21-07-2021", 203, 0, 0, 0, 23, 88, 46, 21, 26),
    @("22-07-2021", 581, 0, 237, 92, 113, 50, 59, 27, 4),
    @("23-07-2021", 171, 0, 0, 0, 39, 16, 38, 43, 36),
    @("26-07-2021", 481, 0, 0, 0, 213, 81, 102, 5, 80),
    @("27-07-2021", 371, 0, 96, 32, 38, 51, 58, 35, 62),
    @("28-07-2021", 917, 96, 321, 31, 124, 25, 201, 47, 73),
    @("29-07-2021", 440, 0, 0, 58, 48, 13, 229, 38, 54),
    @("30-07-2021", 859, 288, 321, 50, 46, 15, 70, 21, 48)
)

$startRow = 126
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")

# 1) Column A holds dates formatted "dd-mm-yyyy" as plain text, matching every
#    existing row. Assigning such a look-alike string straight to .Value makes
#    Excel "smart-convert" it into a real date serial, which would diverge from
#    the source file (shared-string text, no cell style). Writing it as a
#    ="text" formula keeps it as text, then copy/paste-special (values only)
#    bakes the formula result back down to a plain string cell with no style.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $dateText = $data[$i][0]
    $ws.Range("A$row").Formula = "=""$dateText"""
}

$lastRow = $startRow + $data.Count - 1
$ws.Range("A${startRow}:A${lastRow}").Copy()
$ws.Range("A${startRow}:A${lastRow}").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# 2) Columns B:J are the plain numeric totals.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range("$($cols[$c])$row").Value = $rowData[$c + 1]
    }
}
